$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 8 (which holds "extr1"),
# shifting the existing extr1..extr8 rows down by two rows (to rows 10..17).
$ws.Range("A8:E9").EntireRow.Insert()

# Copy the style (bold border, centered) used by column A data cells onto
# the two freshly inserted rows so they match the rest of the table.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(9, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-write the full data block (rows 8-17) with the final values, so the
# index column (A) and data columns are all correct after the insert.
$data = @(
  @(8,  6,  "line7", 14, 11, $true),
  @(9,  7,  "line8", 16, 9,  $true),
  @(10, 8,  "extr1", 5,  12, $false),
  @(11, 9,  "extr2", 5,  9,  $false),
  @(12, 10, "extr3", 10, 11, $false),
  @(13, 11, "extr4", 7,  8,  $false),
  @(14, 12, "extr5", 9,  11, $false),
  @(15, 13, "extr6", 7,  11, $false),
  @(16, 14, "extr7", 5,  7,  $false),
  @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
}
